# Sync attendance_reports, modules_schedules, and assets from main repo - 2026-01-22 09:27:36
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Six "B1D1/B1D2/B1E1/B1E2/B1F1/B1F2" summary rows flip from the stale
#    "Pending" row-style (fill 5 / yellow) to the same "Not Recorded" style
#    (fill 4 / pink) used by the row right above them, and their Status text
#    changes from "Pending" to "Not Recorded".
# ---------------------------------------------------------------------------
$pendingRows = @(184, 211, 238, 265, 292, 319)

foreach ($r in $pendingRows) {
    $srcRow = $r - 1
    $src = $ws.Range("A" + $srcRow + ":I" + $srcRow)
    $src.Copy()
    $dst = $ws.Range("A" + $r + ":I" + $r)
    $dst.PasteSpecial(-4122)   # xlPasteFormats - formatting only, values untouched
    $ws.Range("I" + $r).Value = "Not Recorded"
}
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2) "Recorded By" text for rows recorded by both the automated System and
#    the instructor flips the name order: "System, dnasr281@gmail.com"
#    becomes "dnasr281@gmail.com, System".
# ---------------------------------------------------------------------------
$swapRows = @(8,9,10,12,14,15,17,18,23,34,35,36,38,40,41,43,44,49,60,61,62,64,66,67,69,70,75,86,87,88,90,92,93,95,96,101,112,113,114,116,118,119,121,122,127,138,139,140,142,144,145,147,148,153,164,167,170,174,191,194,197,201,218,221,224,228,245,248,251,255,272,275,278,282,299,302,305,309)

foreach ($r in $swapRows) {
    $ws.Range("G" + $r).Value = "dnasr281@gmail.com, System"
}

# ---------------------------------------------------------------------------
# 3) Class Statistics panel: Missing/Pending session counters updated.
# ---------------------------------------------------------------------------
$ws.Range("L7").Value = 69
$ws.Range("L8").Value = 0

# ---------------------------------------------------------------------------
# 4) Per-group breakdown (rows 21-26): Pending column (P) up by one,
#    Missing column (Q) drops to zero.
# ---------------------------------------------------------------------------
$ws.Range("P21").Value = 7
$ws.Range("Q21").Value = 0
$ws.Range("P22").Value = 7
$ws.Range("Q22").Value = 0
$ws.Range("P23").Value = 7
$ws.Range("Q23").Value = 0
$ws.Range("P24").Value = 8
$ws.Range("Q24").Value = 0
$ws.Range("P25").Value = 7
$ws.Range("Q25").Value = 0
$ws.Range("P26").Value = 7
$ws.Range("Q26").Value = 0
